$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets index 1 / sheet1.xml) - F2, F3, F4
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value = 1364
$wsExhibit.Range("F3").Value = 2896
$wsExhibit.Range("F4").Value = 7

# Sheet "全部类型" (Worksheets index 4 / sheet4.xml) - F3, F4, F5
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F3").Value = 1364
$wsAll.Range("F4").Value = 2896
$wsAll.Range("F5").Value = 7
